# Update course excel file:
# - Split the single "DIANA SCHOOL OF COMMUNITY SERVICES" department label (column C)
#   into per-course-group department names.
# - Simplify the "NSW/QLD/TAS (Currently not accepting enrolments)" location value
#   (column M) into a plain location "NSW/QLD/TAS" plus a separate locationDetail
#   note "Currently not accepting enrolments" (column N) for the Massage rows.
# - Clear out the stale promotion-validity note (column R) for every course row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (department) ---------------------------------------------
$ws.Range("C2:C7").Value   = "Ageing Support"
$ws.Range("C8").Value      = "Community Services"
$ws.Range("C9:C10").Value  = "Early Childhood"
$ws.Range("C11:C12").Value = "Massage"
$ws.Range("C13:C20").Value = "Packages"

# --- Columns M/N (location / locationDetail) for the Massage rows -------
$ws.Range("M11").Value = "NSW/QLD/TAS"
$ws.Range("N11").Value = "Currently not accepting enrolments"

$ws.Range("M12").Value = "NSW/QLD/TAS"
$ws.Range("N12").Value = "Currently not accepting enrolments"

$ws.Range("M20").Value = "NSW/QLD/TAS"
$ws.Range("N20").Value = "Currently not accepting enrolments"

# --- Column R (promotionValidity) ---------------------------------------
$ws.Range("R2:R20").Value = ""
